# Update the canonical URL(s) for this FHIR Implementation Guide StructureDefinition
# export, and refresh the generation Date to match.
#
#   https://hl7.fr/fhir/fr/medication/...   ->   https://hl7.fr/ig/fhir/medication/...

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: URL + Date ---------------------------------------
$metadata = $wb.Worksheets.Item("Metadata")

$metadata.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/StructureDefinition/fr-mp-substance"
$metadata.Range("B8").Value = "2025-05-05T08:11:38+00:00"

# --- "Elements" sheet: Binding Value Set URL ----------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("Z5").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-substance-code"

# The "Binding Value Set" column (Z / column 26) is best-fit width, so it
# needs to be widened slightly (47.09375 -> 47.328125 chars) to accommodate
# the new, slightly longer URL.
$elements.Columns.Item(26).ColumnWidth = 46.43
